$d = $word.ActiveDocument

# Map of paragraph label -> list of strings to append as separate runs
# (the "Fields" paragraph gets its sentence and its final period as two
# distinct runs, matching how the edit was originally authored; every
# other paragraph gets a single run that already ends with a period).
$plan = @(
    @{ Label = "Fields";       Parts = @(" = Used to store the details about object", ".") },
    @{ Label = "Methods";      Parts = @(" = Used to manipulate.") },
    @{ Label = "Constructors"; Parts = @(" = Used to initialize the fields.") },
    @{ Label = "Properties";   Parts = @(" = Used to set or get the values into the private fields.") },
    @{ Label = "Events";       Parts = @(" = Used to raise the notification to other classes.") },
    @{ Label = "Destructors";  Parts = @(" = Used to clear and managed resources.") }
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $label = $text.TrimEnd([char]13, [char]7)

    $entry = $null
    foreach ($candidate in $plan) {
        if ($candidate.Label -eq $label) {
            $entry = $candidate
            break
        }
    }
    if ($entry -eq $null) {
        continue
    }

    # Range covering just the paragraph's visible text (exclude the
    # paragraph mark at the very end).
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End - 1
    $insertionPoint = $d.Range($pEnd, $pEnd)

    # Insert every part in order, remembering the [start,end) of each
    # inserted run so we can separate it from its neighbour afterwards.
    $spans = @()
    foreach ($part in $entry.Parts) {
        $s = $insertionPoint.End
        $insertionPoint.InsertAfter($part)
        $e = $insertionPoint.End
        $spans += , @($s, $e)
        $insertionPoint.Collapse(0)
    }

    # Force each newly inserted span into its own run (instead of being
    # merged into the run before it) by toggling Bold on/off across the
    # span. Do this back-to-front so earlier spans are not disturbed by
    # the Find/Range bookkeeping of later ones.
    for ($j = $spans.Length - 1; $j -ge 0; $j--) {
        $span = $spans[$j]
        $sub = $d.Range($span[0], $span[1])
        $sub.Bold = 1
        $sub.Bold = 0
    }
}
